$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Header text updates (shared strings used by A8 "Volume .. Number .." and
# C9 "Report Covering the Week .. Through ..")
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(19, 2).Text = "32"
$ws.Range("C9").Characters(28, 9).Text = "8/4/2025"
$ws.Range("C9").Characters(49, 8).Text = "8/10/2025"

# ---------------------------------------------------------------------------
# Row 15 - only M15 / N15 change
# ---------------------------------------------------------------------------
$ws.Range("M15").Value = -75
$ws.Range("N15").Value = -83.333333333333

# ---------------------------------------------------------------------------
# Row 16 - C/D/E swap between number and text placeholders, plus new counts
# ---------------------------------------------------------------------------
$ws.Range("D14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("F16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 5
$ws.Range("K16").Value = -82.758620689655
$ws.Range("L16").Value = -50
$ws.Range("M16").Value = -73.684210526315
$ws.Range("N16").Value = -95.689655172413

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 3
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 12
$ws.Range("K17").Value = 20
$ws.Range("L17").Value = 71.428571428571
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = -57.142857142857

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 3
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 5
$ws.Range("H19").Value = 150
$ws.Range("I19").Value = 28
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = -6.666666666666
$ws.Range("L19").Value = -6.666666666666
$ws.Range("M19").Value = -41.666666666666
$ws.Range("N19").Value = -75

# ---------------------------------------------------------------------------
# Row 21 (bold TOTAL row)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 500
$ws.Range("F21").Value = 12
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = 50
$ws.Range("I21").Value = 47
$ws.Range("J21").Value = 71
$ws.Range("K21").Value = -33.802816901408
$ws.Range("L21").Value = -6
$ws.Range("M21").Value = -38.961038961039
$ws.Range("N21").Value = -83.680555555555

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("G24").Value = 3
$ws.Range("J24").Value = 21
$ws.Range("K24").Value = -23.809523809523
$ws.Range("L24").Value = -38.461538461538
$ws.Range("M24").Value = -65.957446808510

# ---------------------------------------------------------------------------
# Row 26 - C/D/E move from text placeholders to real numbers
# ---------------------------------------------------------------------------
$ws.Range("F26").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("F26").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("K15").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("I26").Value = 19
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = -20.833333333333
$ws.Range("L26").Value = -38.709677419354
$ws.Range("M26").Value = 35.714285714285

# ---------------------------------------------------------------------------
# Row 28 - C/D/E move from numbers to text placeholders
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("D14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -50

# ---------------------------------------------------------------------------
# Row 31
# ---------------------------------------------------------------------------
$ws.Range("G31").Value = 1
